# Apply the commit's data changes to the workbook.
#  Sheet "Demand_Projection": add a new row (row 3) with Panama transmission-line demand data.
#  Sheet "Profiles": rescale the existing Costa Rica rows and add new Panama rows for each
#  of the four season/day timeslice combinations (S1D1, S1D2, S2D1, S2D2).
#
# NOTE: this runtime does not persist COM mutations made inside PowerShell `function`
# blocks, so every edit below is written with flat top-level statements / loops.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Demand_Projection -- add Panama row (row 3)
# ---------------------------------------------------------------------------
$wsDemand = $wb.Worksheets.Item("Demand_Projection")

$wsDemand.Cells.Item(3, 1).Value = "Demand"
$wsDemand.Cells.Item(3, 2).Value = "ELCPANXX02"
$wsDemand.Cells.Item(3, 3).Value = "Output demand of transmission lines in Panama"
$wsDemand.Cells.Item(3, 4).Value = "not needed"
$wsDemand.Cells.Item(3, 5).Value = "not needed"
$wsDemand.Cells.Item(3, 6).Value = "not needed"
$wsDemand.Cells.Item(3, 7).Value = "User defined"
$wsDemand.Cells.Item(3, 8).Value = 0

$panamaRow3Values = @(
    62.18,
    66.54000000000001,
    70.90000000000001,
    75.26000000000001,
    79.62,
    84.13,
    88.63,
    93.14,
    97.64,
    102.15,
    106.85,
    111.55,
    116.26,
    120.96,
    125.66,
    130.6,
    135.55,
    140.49,
    145.43,
    150.37,
    155.7,
    161.04,
    166.37,
    171.7,
    177.03,
    182.55,
    188.07,
    193.58,
    199.1,
    204.62
)

for ($i = 0; $i -lt $panamaRow3Values.Count; $i++) {
    $wsDemand.Cells.Item(3, 9 + $i).Value = $panamaRow3Values[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: Profiles -- update existing rows & append new Panama rows
# ---------------------------------------------------------------------------
$wsProfiles = $wb.Worksheets.Item("Profiles")

# Metadata (columns A-I) for rows 2-9, in order.
$profileMeta = @(
    @{ Row = 2; Timeslice = "S1D1"; FuelTech = "ELCCRIXX02"; Name = "Output demand of transmission lines in Costa Rica"; Value = 0.23 },
    @{ Row = 3; Timeslice = "S1D1"; FuelTech = "ELCPANXX02"; Name = "Output demand of transmission lines in Panama";     Value = 0.21 },
    @{ Row = 4; Timeslice = "S1D2"; FuelTech = "ELCCRIXX02"; Name = "Output demand of transmission lines in Costa Rica"; Value = 0.27 },
    @{ Row = 5; Timeslice = "S1D2"; FuelTech = "ELCPANXX02"; Name = "Output demand of transmission lines in Panama";     Value = 0.27 },
    @{ Row = 6; Timeslice = "S2D1"; FuelTech = "ELCCRIXX02"; Name = "Output demand of transmission lines in Costa Rica"; Value = 0.23 },
    @{ Row = 7; Timeslice = "S2D1"; FuelTech = "ELCPANXX02"; Name = "Output demand of transmission lines in Panama";     Value = 0.23 },
    @{ Row = 8; Timeslice = "S2D2"; FuelTech = "ELCCRIXX02"; Name = "Output demand of transmission lines in Costa Rica"; Value = 0.27 },
    @{ Row = 9; Timeslice = "S2D2"; FuelTech = "ELCPANXX02"; Name = "Output demand of transmission lines in Panama";     Value = 0.29 }
)

foreach ($entry in $profileMeta) {
    $row = $entry.Row

    $wsProfiles.Cells.Item($row, 1).Value = $entry.Timeslice
    $wsProfiles.Cells.Item($row, 2).Value = "Demand"
    $wsProfiles.Cells.Item($row, 3).Value = $entry.FuelTech
    $wsProfiles.Cells.Item($row, 4).Value = $entry.Name
    $wsProfiles.Cells.Item($row, 5).Value = "not needed"
    $wsProfiles.Cells.Item($row, 6).Value = "not needed"
    $wsProfiles.Cells.Item($row, 7).Value = "not needed"
    $wsProfiles.Cells.Item($row, 8).Value = "User defined"
    $wsProfiles.Cells.Item($row, 9).Value = 0

    for ($col = 10; $col -le 39; $col++) {
        $wsProfiles.Cells.Item($row, $col).Value = $entry.Value
    }
}

Write-Output "Done applying edits."
